# Apply updated cryptocurrency price/volume data per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.315.79"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "3.071.13"
$ws.Range("E3").Value = "  +3.16%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'579.29"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "'167.94"
$ws.Range("E6").Value = "  +4.15%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.068.09"
$ws.Range("E8").Value = "  +3.10%  "

$ws.Range("E9").Value = "  +1.38%  "

$ws.Range("D10").Value = "'6.65"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "'0.485"
$ws.Range("E12").Value = "  +6.95%  "

$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").Value = "'36.67"
$ws.Range("E14").Value = "  +6.58%  "

$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "3.584.70"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("D17").Value = "66.323.67"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "'7.17"
$ws.Range("E18").Value = "  +4.27%  "

$ws.Range("D19").Value = "3.075.63"
$ws.Range("E19").Value = "  +3.43%  "

$ws.Range("D20").Value = "'16.21"
$ws.Range("E20").Value = "  +17.61%  "

$ws.Range("D21").Value = "'463.35"
$ws.Range("E21").Value = "  +2.64%  "

$ws.Range("D22").Value = "'0.711"
$ws.Range("E22").Value = "  +4.69%  "

$ws.Range("D23").Value = "'7.44"
$ws.Range("E23").Value = "  +2.23%  "

$ws.Range("D24").Value = "'83.20"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("D25").Value = "'12.82"
$ws.Range("E25").Value = "  +5.26%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'8.07"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").Value = "'2.41"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("E31").Value = "  +3.17%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'28.25"
$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("E34").Value = "  +5.13%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.87"
$ws.Range("E36").Value = "  +1.83%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.993"
$ws.Range("E37").Value = "  +0.88%  "

$ws.Range("D38").Value = "'48.22"
$ws.Range("E38").Value = "  +10.09%  "

$ws.Range("D39").Value = "'49.95"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.04"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.314"
$ws.Range("E41").Value = "  +4.88%  "

$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  +1.59%  "

$ws.Range("D44").Value = "'8.63"
$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("E45").Value = "  +1.85%  "

$ws.Range("D46").Value = "'381.32"
$ws.Range("E46").Value = "  -2.17%  "

$ws.Range("D47").Value = "2.763.79"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").Value = "'134.09"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "'24.46"
$ws.Range("E50").Value = "  +5.44%  "

$ws.Range("E51").Value = "  +4.02%  "
